{"js": "// Apply the Historias_De_Usuario_V3.0 wording updates.\n// Each edit is a precise text replacement located with body.search(),\n// then rewritten in place with Range.insertText(text, \"Replace\").\n\nconst body = context.document.body;\n\nasync function replaceOnce(oldText, newText) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\n      \"expected exactly 1 match for \" + JSON.stringify(oldText) + \" but found \" + results.items.length\n    );\n  }\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n\n// 1) \"Quiero poder convertir...\" -> \"Quiero convertir...\"\nawait replaceOnce(\n  \" poder convertir un n\u00famero a su equivalente en palabras en kichwa \",\n  \" convertir un n\u00famero a su equivalente en palabras en kichwa \"\n);\n\n// 2) \"...entrada, cuando se llama a la funci\u00f3n...\" -> \"...entrada, al momento de llamar a la funci\u00f3n...\"\nawait replaceOnce(\n  \"cuando se llama a la funci\u00f3n de conversi\u00f3n a kichwa con el n\u00famero\",\n  \"al momento de llamar a la funci\u00f3n de conversi\u00f3n a kichwa con el n\u00famero\"\n);\n\n// 3) \"sin importar si es negativo opositivo, se detecta que es un n\u00famero...\" ->\n//    \"sin importar si es negativo o positivo, detecte que el valor ingresado es un n\u00famero...\"\nawait replaceOnce(\n  \"sin importar si es negativo opositivo, se detecta que es un n\u00famero y no otro tipo de dato\",\n  \"sin importar si es negativo o positivo, detecte que el valor ingresado es un n\u00famero y no otro tipo de dato\"\n);\n\n// 4) \"...no proporciona un n\u00famero como entrada a la funci\u00f3n...\" -> \"...como entrada en la funci\u00f3n...\"\n// (the sibling story about the 0-9999 range has near-identical wording, so anchor on the\n// full sentence that is unique to this paragraph)\nawait replaceOnce(\n  \"Verificar que cuando un usuario no proporciona un n\u00famero como entrada a la funci\u00f3n de conversi\u00f3n a kichwa, se genere un mensaje de error indicando que se requiere un n\u00famero v\u00e1lido.\",\n  \"Verificar que cuando un usuario no proporciona un n\u00famero como entrada en la funci\u00f3n de conversi\u00f3n a kichwa, se genere un mensaje de error indicando que se requiere un n\u00famero v\u00e1lido.\"\n);\n\n// 5) \"...indique como usar el programa...\" -> \"...indique el funcionamiento del programa...\"\nawait replaceOnce(\n  \"El programa debe venir con un manual de usuario que indique como usar el programa. (Capacidad para reconocer su adecuaci\u00f3n)\",\n  \"El programa debe venir con un manual de usuario que indique el funcionamiento del programa. (Capacidad para reconocer su adecuaci\u00f3n)\"\n);\n\n// 6) \"...puede manejar n\u00fameros...\" -> \"...es capaz de manejar n\u00fameros...\"\nawait replaceOnce(\n  \"Comprobar que la funci\u00f3n de conversi\u00f3n a kichwa puede manejar n\u00fameros \",\n  \"Comprobar que la funci\u00f3n de conversi\u00f3n a kichwa es capaz de manejar n\u00fameros \"\n);\n\n// 7) \"...no consuma una excesiva cantidad de recursos del procesador que provoque la ralentizaci\u00f3n del dispositivo...\" ->\n//    \"...no consuma una cantidad excesiva de recursos del procesador que provoque una disminuci\u00f3n en el rendimiento del dispositivo...\"\nawait replaceOnce(\n  \"no consuma una excesiva cantidad de recursos del procesador que provoque la ralentizaci\u00f3n del dispositivo. (Utilizaci\u00f3n de recursos)\",\n  \"no consuma una cantidad excesiva de recursos del procesador que provoque una disminuci\u00f3n en el rendimiento del dispositivo. (Utilizaci\u00f3n de recursos)\"\n);\n\n// 8) \"Los mensajes se presentan de forma legible y bien formateada en la consola, independientemente.\" ->\n//    \"Los mensajes se presentan de forma legible y ordenada en la consola.\"\nawait replaceOnce(\n  \"Los mensajes se presentan de forma legible y bien formateada en la consola, independientemente.\",\n  \"Los mensajes se presentan de forma legible y ordenada en la consola.\"\n);\n", "ps1": "# Apply the Historias_De_Usuario_V3.0 wording updates via Word COM Find/Replace.\n# Each call locates one exact (and unique) run of text in the document body and\n# replaces it in place, leaving formatting/runs around it untouched.\n\n$d = $word.ActiveDocument\n\nfunction Replace-ExactText {\n    param(\n        [string]$FindText,\n        [string]$ReplaceText\n    )\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $result = $find.Execute(\n        $FindText,   # FindText\n        $false,      # MatchCase\n        $false,      # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        1,           # Wrap (wdFindContinue)\n        $false,      # Format\n        $ReplaceText,# ReplaceWith\n        2            # Replace (wdReplaceAll)\n    )\n\n    if (-not $result) {\n        throw \"Could not find text: $FindText\"\n    }\n}\n\n# 1) \"Quiero poder convertir...\" -> \"Quiero convertir...\"\nReplace-ExactText \" poder convertir un n\u00famero a su equivalente en palabras en kichwa \" \" convertir un n\u00famero a su equivalente en palabras en kichwa \"\n\n# 2) \"...entrada, cuando se llama a la funci\u00f3n...\" -> \"...entrada, al momento de llamar a la funci\u00f3n...\"\nReplace-ExactText \"cuando se llama a la funci\u00f3n de conversi\u00f3n a kichwa con el n\u00famero\" \"al momento de llamar a la funci\u00f3n de conversi\u00f3n a kichwa con el n\u00famero\"\n\n# 3) \"sin importar si es negativo opositivo, se detecta que es un n\u00famero...\" ->\n#    \"sin importar si es negativo o positivo, detecte que el valor ingresado es un n\u00famero...\"\nReplace-ExactText \"sin importar si es negativo opositivo, se detecta que es un n\u00famero y no otro tipo de dato\" \"sin importar si es negativo o positivo, detecte que el valor ingresado es un n\u00famero y no otro tipo de dato\"\n\n# 4) \"...no proporciona un n\u00famero como entrada a la funci\u00f3n...\" -> \"...como entrada en la funci\u00f3n...\"\n# (the sibling story about the 0-9999 range has near-identical wording, so anchor on the\n# full sentence that is unique to this paragraph)\nReplace-ExactText \"Verificar que cuando un usuario no proporciona un n\u00famero como entrada a la funci\u00f3n de conversi\u00f3n a kichwa, se genere un mensaje de error indicando que se requiere un n\u00famero v\u00e1lido.\" \"Verificar que cuando un usuario no proporciona un n\u00famero como entrada en la funci\u00f3n de conversi\u00f3n a kichwa, se genere un mensaje de error indicando que se requiere un n\u00famero v\u00e1lido.\"\n\n# 5) \"...indique como usar el programa...\" -> \"...indique el funcionamiento del programa...\"\nReplace-ExactText \"El programa debe venir con un manual de usuario que indique como usar el programa. (Capacidad para reconocer su adecuaci\u00f3n)\" \"El programa debe venir con un manual de usuario que indique el funcionamiento del programa. (Capacidad para reconocer su adecuaci\u00f3n)\"\n\n# 6) \"...puede manejar n\u00fameros...\" -> \"...es capaz de manejar n\u00fameros...\"\nReplace-ExactText \"Comprobar que la funci\u00f3n de conversi\u00f3n a kichwa puede manejar n\u00fameros\" \"Comprobar que la funci\u00f3n de conversi\u00f3n a kichwa es capaz de manejar n\u00fameros\"\n\n# 7) \"...no consuma una excesiva cantidad de recursos del procesador que provoque la ralentizaci\u00f3n del dispositivo...\" ->\n#    \"...no consuma una cantidad excesiva de recursos del procesador que provoque una disminuci\u00f3n en el rendimiento del dispositivo...\"\nReplace-ExactText \"no consuma una excesiva cantidad de recursos del procesador que provoque la ralentizaci\u00f3n del dispositivo. (Utilizaci\u00f3n de recursos)\" \"no consuma una cantidad excesiva de recursos del procesador que provoque una disminuci\u00f3n en el rendimiento del dispositivo. (Utilizaci\u00f3n de recursos)\"\n\n# 8) \"Los mensajes se presentan de forma legible y bien formateada en la consola, independientemente.\" ->\n#    \"Los mensajes se presentan de forma legible y ordenada en la consola.\"\nReplace-ExactText \"Los mensajes se presentan de forma legible y bien formateada en la consola, independientemente.\" \"Los mensajes se presentan de forma legible y ordenada en la consola.\"\n"}
